$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("F2").Value = 91
$ws.Range("J2").Value = 120
$ws.Range("L2").Value = 127
$ws.Range("D3").Value = 132
$ws.Range("F3").Value = 135
$ws.Range("J3").Value = 226
$ws.Range("K3").Value = 216
$ws.Range("F4").Value = 8
$ws.Range("B6").Value = 370
$ws.Range("C6").Value = 472
$ws.Range("D6").Value = 410
$ws.Range("E6").Value = 464
$ws.Range("F6").Value = 520
$ws.Range("G6").Value = 433
$ws.Range("H6").Value = 436
$ws.Range("I6").Value = 497
$ws.Range("J6").Value = 411
$ws.Range("K6").Value = 499
$ws.Range("L6").Value = 424
$ws.Range("B7").Value = 495
$ws.Range("C7").Value = 626
$ws.Range("D7").Value = 639
$ws.Range("E7").Value = 687
$ws.Range("F7").Value = 755
$ws.Range("G7").Value = 662
$ws.Range("H7").Value = 708
$ws.Range("I7").Value = 829
$ws.Range("J7").Value = 780
$ws.Range("K7").Value = 881
$ws.Range("L7").Value = 817

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("C6").Value = 34
$ws.Range("E6").Value = 52
$ws.Range("I6").Value = 32
$ws.Range("K6").Value = 25
$ws.Range("C7").Value = 39
$ws.Range("E7").Value = 65
$ws.Range("I7").Value = 48
$ws.Range("K7").Value = 45
$ws.Range("K3").Value = 17

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("D6").Value = 21
$ws.Range("G6").Value = 15
$ws.Range("D7").Value = 36
$ws.Range("G7").Value = 27
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 14

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 15

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("D3").Value = 17
$ws.Range("H6").Value = 26
$ws.Range("D7").Value = 45
$ws.Range("H7").Value = 45

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("D4").Value = 6
$ws.Range("D5").Value = 10

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("F5").Value = 14
$ws.Range("B8").Value = 30
$ws.Range("F8").Value = 50
$ws.Range("L8").Value = 31
$ws.Range("D10").Value = 3
$ws.Range("B16").Value = 2
$ws.Range("D19").Value = 27
$ws.Range("F19").Value = 23
$ws.Range("D28").Value = 45
$ws.Range("H28").Value = 45
$ws.Range("F29").Value = 13
$ws.Range("C32").Value = 39
$ws.Range("E32").Value = 65
$ws.Range("I32").Value = 48
$ws.Range("K32").Value = 45
$ws.Range("D36").Value = 36
$ws.Range("G36").Value = 27
$ws.Range("K36").Value = 67
$ws.Range("J42").Value = 13
$ws.Range("D47").Value = 14
$ws.Range("F47").Value = 17
$ws.Range("J47").Value = 16
$ws.Range("F51").Value = 7
$ws.Range("E53").Value = 81
$ws.Range("F53").Value = 80
$ws.Range("I53").Value = 124
$ws.Range("J53").Value = 120
$ws.Range("J54").Value = 10
$ws.Range("D61").Value = 3
$ws.Range("L61").Value = 1
$ws.Range("E62").Value = 7
$ws.Range("C63").Value = 7
$ws.Range("C65").Value = 22
$ws.Range("F70").Value = 24
$ws.Range("I75").Value = 2
$ws.Range("F76").Value = 19
$ws.Range("K76").Value = 29
$ws.Range("G77").Value = 24
$ws.Range("B80").Value = 15
$ws.Range("D82").Value = 10
$ws.Range("C85").Value = 15
$ws.Range("J91").Value = 7
$ws.Range("G97").Value = 5
$ws.Range("B98").Value = 495
$ws.Range("C98").Value = 626
$ws.Range("D98").Value = 639
$ws.Range("E98").Value = 687
$ws.Range("F98").Value = 755
$ws.Range("G98").Value = 662
$ws.Range("H98").Value = 708
$ws.Range("I98").Value = 829
$ws.Range("J98").Value = 780
$ws.Range("K98").Value = 881
$ws.Range("L98").Value = 817

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("F2").Value = 7
$ws.Range("E6").Value = 63
$ws.Range("F6").Value = 59
$ws.Range("I6").Value = 79
$ws.Range("J6").Value = 60
$ws.Range("E7").Value = 81
$ws.Range("F7").Value = 80
$ws.Range("I7").Value = 124
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 7
$ws.Range("F4").Value = 2

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K5").Value = 14
$ws.Range("F6").Value = 19
$ws.Range("K6").Value = 29

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("C5").Value = 18
$ws.Range("C6").Value = 22

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("E6").Value = 5

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("E7").Value = 7

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("F5").Value = 12
$ws.Range("F6").Value = 13

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("G6").Value = 14
$ws.Range("G7").Value = 24

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("H4").Value = 11
$ws.Range("H5").Value = 13

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 3

$ws = $wb.Worksheets.Item('New City')
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 7
$ws.Range("G5").Value = 3

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("G6").Value = 5

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("D6").Value = 15
$ws.Range("F6").Value = 16
$ws.Range("D7").Value = 27
$ws.Range("F7").Value = 23

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 7

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("C4").Value = 12
$ws.Range("C5").Value = 15

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J3").Value = 2

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J5").Value = 10

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 11
$ws.Range("D5").Value = 12
$ws.Range("F5").Value = 10
$ws.Range("D6").Value = 14
$ws.Range("F6").Value = 17
$ws.Range("J6").Value = 16

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("F5").Value = 18
$ws.Range("F6").Value = 24

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 2

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 4
$ws.Range("F3").Value = 6
$ws.Range("B6").Value = 21
$ws.Range("B7").Value = 30
$ws.Range("F7").Value = 50
$ws.Range("L7").Value = 31
